$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3385.7144
$ws.Range("I64").Value = 3183.3333
$ws.Range("J64").Value = 3537.5
$ws.Range("K64").Value = 3183.3333
$ws.Range("L64").Value = 3537.5
$ws.Range("M64").Value = -2935.3333
$ws.Range("N64").Value = -4033.5
# Row 67
$ws.Range("H67").Value = 3385.7144
$ws.Range("I67").Value = 3183.3333
$ws.Range("J67").Value = 3537.5
$ws.Range("K67").Value = 3183.3333
$ws.Range("L67").Value = 3537.5
$ws.Range("M67").Value = -2325.3333
$ws.Range("N67").Value = -5253.5
# Row 74
$ws.Range("H74").Value = 3032.9788
$ws.Range("I74").Value = 2963.3157
$ws.Range("J74").Value = 3327.111
$ws.Range("K74").Value = 2963.3157
$ws.Range("L74").Value = 3327.111
$ws.Range("M74").Value = -2027.3157
$ws.Range("N74").Value = -5199.111
# Row 76
$ws.Range("H76").Value = 3065.9
$ws.Range("I76").Value = 3057.375
$ws.Range("K76").Value = 3057.375
$ws.Range("M76").Value = -2742.375
# Row 77
$ws.Range("H77").Value = 3032.9788
$ws.Range("I77").Value = 2963.3157
$ws.Range("J77").Value = 3327.111
$ws.Range("K77").Value = 14816.5785
$ws.Range("L77").Value = 16635.555
$ws.Range("M77").Value = -10136.5785
$ws.Range("N77").Value = -25995.555
# Row 79
$ws.Range("H79").Value = 3065.9
$ws.Range("I79").Value = 3057.375
$ws.Range("K79").Value = 3057.375
$ws.Range("M79").Value = -1965.375
# Row 132
$ws.Range("H132").Value = 4503.8125
$ws.Range("I132").Value = 1367.8077
$ws.Range("J132").Value = 18093.166
$ws.Range("K132").Value = 4103.4231
$ws.Range("L132").Value = 54279.49800000001
$ws.Range("M132").Value = -1573.4231
$ws.Range("N132").Value = -59339.49800000001
# Row 135
$ws.Range("H135").Value = 250002770
$ws.Range("I135").Value = 2060
$ws.Range("J135").Value = 333336320
$ws.Range("K135").Value = 18540
$ws.Range("L135").Value = 3000026880
$ws.Range("M135").Value = -16005
$ws.Range("N135").Value = -3000031950

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3516.76
$ws.Range("I32").Value = 3257.0417
$ws.Range("J32").Value = 9750
$ws.Range("K32").Value = 3257.0417
$ws.Range("L32").Value = 9750
$ws.Range("M32").Value = -2970.0417
$ws.Range("N32").Value = -10324
# Row 74
$ws.Range("H74").Value = 201320.8
$ws.Range("I74").Value = 251456
$ws.Range("J74").Value = 780
$ws.Range("K74").Value = 251456
$ws.Range("L74").Value = 780
$ws.Range("M74").Value = -250582
$ws.Range("N74").Value = -2528
# Row 77
$ws.Range("H77").Value = 201320.8
$ws.Range("I77").Value = 251456
$ws.Range("J77").Value = 780
$ws.Range("K77").Value = 1257280
$ws.Range("L77").Value = 3900
$ws.Range("M77").Value = -1252912
$ws.Range("N77").Value = -12636
# Row 98
$ws.Range("H98").Value = 24450
$ws.Range("J98").Value = 24450
$ws.Range("L98").Value = 24450
$ws.Range("N98").Value = -30440
# Row 102
$ws.Range("H102").Value = 1146
$ws.Range("I102").Value = 1146
$ws.Range("K102").Value = 1146
$ws.Range("M102").Value = 476

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 27666.666
$ws.Range("J100").Value = 27666.666
$ws.Range("L100").Value = 27666.666
$ws.Range("N100").Value = -29830.666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 49.5
$ws.Range("I7").Value = 40.076923
$ws.Range("J7").Value = 60.636364
$ws.Range("K7").Value = 40.076923
$ws.Range("L7").Value = 60.636364
$ws.Range("M7").Value = 72.92307700000001
$ws.Range("N7").Value = -286.636364
# Row 22
$ws.Range("H22").Value = 695.4545000000001
$ws.Range("I22").Value = 575
$ws.Range("J22").Value = 764.2857
$ws.Range("K22").Value = 575
$ws.Range("L22").Value = 764.2857
$ws.Range("M22").Value = -225
$ws.Range("N22").Value = -1464.2857
# Row 94
$ws.Range("H94").Value = 4103.75
$ws.Range("I94").Value = 5762
$ws.Range("J94").Value = 2445.5
$ws.Range("K94").Value = 5762
$ws.Range("L94").Value = 2445.5
$ws.Range("M94").Value = -5311
$ws.Range("N94").Value = -3347.5
# Row 100
$ws.Range("H100").Value = 31538.092
$ws.Range("J100").Value = 31538.092
$ws.Range("L100").Value = 31538.092
$ws.Range("N100").Value = -33702.092

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 472.34375
$ws.Range("I5").Value = 297.08334
$ws.Range("K5").Value = 891.2500200000001
$ws.Range("M5").Value = -779.2500200000001
# Row 97
$ws.Range("H97").Value = 1105.5454
$ws.Range("I97").Value = 498.69232
$ws.Range("K97").Value = 1496.07696
$ws.Range("M97").Value = -1000.07696
# Row 122
$ws.Range("H122").Value = 24447.262
$ws.Range("I122").Value = 32745.766
$ws.Range("J122").Value = 934.8333
$ws.Range("K122").Value = 294711.894
$ws.Range("L122").Value = 8413.4997
$ws.Range("M122").Value = -292261.894
$ws.Range("N122").Value = -13313.4997
# Row 135
$ws.Range("H135").Value = 472.34375
$ws.Range("I135").Value = 297.08334
$ws.Range("K135").Value = 2673.75006
$ws.Range("M135").Value = -138.7500600000003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 2413.1667
$ws.Range("I31").Value = 2413.1667
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2413.1667
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2121.1667
$ws.Range("N31").ClearContents()
# Row 37
$ws.Range("H37").Value = 2413.1667
$ws.Range("I37").Value = 2413.1667
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2413.1667
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2136.1667
$ws.Range("N37").ClearContents()
# Row 132
$ws.Range("H132").Value = 2862742
$ws.Range("I132").Value = 4062.9412
$ws.Range("J132").Value = 5562605.5
$ws.Range("K132").Value = 12188.8236
$ws.Range("L132").Value = 16687816.5
$ws.Range("M132").Value = -9658.8236
$ws.Range("N132").Value = -16692876.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 35
$ws.Range("H35").Value = 12687.333
$ws.Range("I35").Value = 1531
$ws.Range("K35").Value = 1531
$ws.Range("M35").Value = -1195
# Row 61
$ws.Range("H61").Value = 2138.5
$ws.Range("I61").Value = 1809.1428
$ws.Range("J61").Value = 2599.6
$ws.Range("K61").Value = 1809.1428
$ws.Range("L61").Value = 2599.6
$ws.Range("M61").Value = -1607.1428
$ws.Range("N61").Value = -3003.6
# Row 111
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180
# Row 113
$ws.Range("H113").Value = 2138.5
$ws.Range("I113").Value = 1809.1428
$ws.Range("J113").Value = 2599.6
$ws.Range("K113").Value = 1809.1428
$ws.Range("L113").Value = 2599.6
$ws.Range("M113").Value = 360.8571999999999
$ws.Range("N113").Value = -6939.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 97
$ws.Range("H97").Value = 45457.332
$ws.Range("J97").Value = 45457.332
$ws.Range("L97").Value = 45457.332
$ws.Range("N97").Value = -47439.332
